# The "客単価" (avg spend per customer) column H was being computed without
# dividing by the customer count (column E, count_客構成), so every value in
# H was inflated by a factor of E. This corrects H by dividing each numeric
# H value by the corresponding E value on the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC分析_客構成")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    $eCell = $ws.Cells.Item($r, 5)

    $hVal = $hCell.Value()
    $eVal = $eCell.Value()

    if ($eVal -ne $null -and $eVal -ne 0) {
        if ($hVal -is [double] -or $hVal -is [int]) {
            $hCell.Value = $hVal / $eVal
        }
    }
}
